$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3647.6296
$ws.Range("I40").Value = 2044
$ws.Range("J40").Value = 4322.8423
$ws.Range("K40").Value = 2044
$ws.Range("L40").Value = 4322.8423
$ws.Range("M40").Value = -1869
$ws.Range("N40").Value = -4672.8423
$ws.Range("H58").Value = 1242
$ws.Range("J58").Value = 2429
$ws.Range("L58").Value = 7287
$ws.Range("N58").Value = -7587
$ws.Range("H81").Value = 32500
$ws.Range("I81").Value = 15000
$ws.Range("J81").Value = 38333.332
$ws.Range("K81").Value = 15000
$ws.Range("L81").Value = 38333.332
$ws.Range("M81").Value = -14002
$ws.Range("N81").Value = -40329.332
$ws.Range("H84").Value = 32500
$ws.Range("I84").Value = 15000
$ws.Range("J84").Value = 38333.332
$ws.Range("K84").Value = 45000
$ws.Range("L84").Value = 114999.996
$ws.Range("M84").Value = -40008
$ws.Range("N84").Value = -124983.996
$ws.Range("H132").Value = 1928.7705
$ws.Range("I132").Value = 1909.0892
$ws.Range("K132").Value = 5727.267599999999
$ws.Range("M132").Value = -3197.267599999999
$ws.Range("H137").Value = 60845.934
$ws.Range("I137").Value = 91561.45
$ws.Range("K137").Value = 274684.35
$ws.Range("M137").Value = -272134.35

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6546.99
$ws.Range("I32").Value = 4997.3013
$ws.Range("J32").Value = 27135.715
$ws.Range("K32").Value = 4997.3013
$ws.Range("L32").Value = 27135.715
$ws.Range("M32").Value = -4710.3013
$ws.Range("N32").Value = -27709.715
$ws.Range("H45").Value = 5717181.5
$ws.Range("I45").Value = 8929628
$ws.Range("K45").Value = 8929628
$ws.Range("M45").Value = -8929251
$ws.Range("H61").Value = 5631.7144
$ws.Range("I61").Value = 5237
$ws.Range("K61").Value = 5237
$ws.Range("M61").Value = -5025
$ws.Range("H74").Value = 32801.7
$ws.Range("I74").Value = 2599.4736
$ws.Range("K74").Value = 2599.4736
$ws.Range("M74").Value = -1725.4736
$ws.Range("H77").Value = 32801.7
$ws.Range("I77").Value = 2599.4736
$ws.Range("K77").Value = 12997.368
$ws.Range("M77").Value = -8629.367999999999
$ws.Range("H122").Value = 1471931.1
$ws.Range("I122").Value = 1317385.9
$ws.Range("K122").Value = 3952157.7
$ws.Range("M122").Value = -3949707.7
$ws.Range("H132").Value = 47384.1
$ws.Range("I132").Value = 2058.0715
$ws.Range("J132").Value = 153144.83
$ws.Range("K132").Value = 6174.2145
$ws.Range("L132").Value = 459434.49
$ws.Range("M132").Value = -3644.2145
$ws.Range("N132").Value = -464494.49
$ws.Range("H136").Value = 5631.7144
$ws.Range("I136").Value = 5237
$ws.Range("K136").Value = 15711
$ws.Range("M136").Value = -13161
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2152.1892
$ws.Range("I20").Value = 1613.12
$ws.Range("J20").Value = 3275.25
$ws.Range("K20").Value = 1613.12
$ws.Range("L20").Value = 3275.25
$ws.Range("M20").Value = -1366.12
$ws.Range("N20").Value = -3769.25
$ws.Range("H86").Value = 15479994
$ws.Range("J86").Value = 4526.8237
$ws.Range("L86").Value = 4526.8237
$ws.Range("N86").Value = -6772.8237
$ws.Range("H89").Value = 15479994
$ws.Range("J89").Value = 4526.8237
$ws.Range("L89").Value = 22634.1185
$ws.Range("N89").Value = -33866.1185
$ws.Range("H133").Value = 158997.5
$ws.Range("J133").Value = 158997.5
$ws.Range("L133").Value = 158997.5
$ws.Range("N133").Value = -169117.5
$ws.Range("H139").Value = 179995
$ws.Range("J139").Value = 179995
$ws.Range("L139").Value = 179995
$ws.Range("N139").Value = -190275

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18224.951
$ws.Range("I31").Value = 2448.25
$ws.Range("J31").Value = 56790.223
$ws.Range("K31").Value = 2448.25
$ws.Range("L31").Value = 56790.223
$ws.Range("M31").Value = -2153.25
$ws.Range("N31").Value = -57380.223
$ws.Range("H34").Value = 18224.951
$ws.Range("I34").Value = 2448.25
$ws.Range("J34").Value = 56790.223
$ws.Range("K34").Value = 2448.25
$ws.Range("L34").Value = 56790.223
$ws.Range("M34").Value = -2246.25
$ws.Range("N34").Value = -57194.223
$ws.Range("H58").Value = 5803.8965
$ws.Range("I58").Value = 6740.7896
$ws.Range("K58").Value = 6740.7896
$ws.Range("M58").Value = -6537.7896
$ws.Range("H99").Value = 4176.467
$ws.Range("I99").Value = 3785.6365
$ws.Range("K99").Value = 3785.6365
$ws.Range("M99").Value = -2287.6365
$ws.Range("H126").Value = 4176.467
$ws.Range("I126").Value = 3785.6365
$ws.Range("K126").Value = 11356.9095
$ws.Range("M126").Value = -8886.9095
$ws.Range("H136").Value = 5803.8965
$ws.Range("I136").Value = 6740.7896
$ws.Range("K136").Value = 20222.3688
$ws.Range("M136").Value = -17672.3688
$ws.Range("H141").Value = 132999.8
$ws.Range("J141").Value = 132999.8
$ws.Range("L141").Value = 132999.8
$ws.Range("N141").Value = -143359.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3882.182
$ws.Range("I3").Value = 3764.5
$ws.Range("J3").Value = 3908.3333
$ws.Range("K3").Value = 11293.5
$ws.Range("L3").Value = 11724.9999
$ws.Range("M3").Value = -11181.5
$ws.Range("N3").Value = -11948.9999
$ws.Range("H18").Value = 1041.8889
$ws.Range("I18").Value = 255.4
$ws.Range("K18").Value = 766.2
$ws.Range("M18").Value = -597.2
$ws.Range("H51").Value = 5360.643
$ws.Range("I51").Value = 1200
$ws.Range("K51").Value = 3600
$ws.Range("M51").Value = -3140
$ws.Range("H55").Value = 83337460
$ws.Range("J55").Value = 95399.09
$ws.Range("L55").Value = 286197.27
$ws.Range("N55").Value = -286551.27
$ws.Range("H134").Value = 499.5
$ws.Range("I134").Value = 499.5
$ws.Range("K134").Value = 1498.5
$ws.Range("M134").Value = 3571.5
$ws.Range("H136").Value = 6860.5557
$ws.Range("I136").Value = 6468.125
$ws.Range("K136").Value = 19404.375
$ws.Range("M136").Value = -14304.375
$ws.Range("H137").Value = 7825.154
$ws.Range("J137").Value = 10062.375
$ws.Range("L137").Value = 30187.125
$ws.Range("N137").Value = -40387.125
$ws.Range("H140").Value = 1729.1
$ws.Range("I140").Value = 930.125
$ws.Range("J140").Value = 4925
$ws.Range("K140").Value = 2790.375
$ws.Range("L140").Value = 14775
$ws.Range("M140").Value = 2389.625
$ws.Range("N140").Value = -25135

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 810768
$ws.Range("J80").Value = 182325.84
$ws.Range("L80").Value = 182325.84
$ws.Range("N80").Value = -184321.84
$ws.Range("H83").Value = 810768
$ws.Range("J83").Value = 182325.84
$ws.Range("L83").Value = 911629.2
$ws.Range("N83").Value = -921613.2
$ws.Range("H132").Value = 3423.653
$ws.Range("I132").Value = 3282.5144
$ws.Range("K132").Value = 9847.5432
$ws.Range("M132").Value = -7317.5432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9511.241
$ws.Range("I132").Value = 10372.375
$ws.Range("J132").Value = 5377.8
$ws.Range("K132").Value = 31117.125
$ws.Range("L132").Value = 16133.4
$ws.Range("M132").Value = -28587.125
$ws.Range("N132").Value = -21193.4
$ws.Range("H136").Value = 225867
$ws.Range("I136").Value = 253100.38
$ws.Range("K136").Value = 759301.14
$ws.Range("M136").Value = -756751.14
$ws.Range("H139").Value = 93972
$ws.Range("J139").Value = 93972
$ws.Range("L139").Value = 93972
$ws.Range("N139").Value = -104252

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 16851366
$ws.Range("I132").Value = 18522002
$ws.Range("J132").Value = 1815647.9
$ws.Range("K132").Value = 55566006
$ws.Range("L132").Value = 5446943.699999999
$ws.Range("M132").Value = -55563476
$ws.Range("N132").Value = -5452003.699999999
$ws.Range("H139").Value = 70376.39999999999
$ws.Range("J139").Value = 70376.39999999999
$ws.Range("L139").Value = 70376.39999999999
$ws.Range("N139").Value = -80656.39999999999
